$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "name "
$ws.Range("B1").Value = "fName"
$ws.Range("C1").Value = "cource"
$ws.Range("A2").Value = "kundan "
$ws.Range("B2").Value = "vinod  "
$ws.Range("C2").Value = "dsa"

$ws.Range("C2").Select()
